# Weekly update: insert two new price rows (week of 2021-09-22 / serial 44461)
# at the top of the Piña series, pushing the previously-first rows (131-136)
# down to become rows 133-138.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 131..136 down by two, preserving all their data/formatting.
$ws.Rows.Item(131).Insert()
$ws.Rows.Item(131).Insert()

# New row 131 - Primera, $/caja 12 unidades
$ws.Range("A131").Value = 7
$ws.Range("B131").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C131").Value = "Ñuble"
$ws.Range("D131").Value = 44461
$ws.Range("E131").Value = 16
$ws.Range("F131").Value = "Fruta"
$ws.Range("G131").Value = 100108
$ws.Range("H131").Value = "Tropicales y subtropicales"
$ws.Range("I131").Value = 100108005
$ws.Range("J131").Value = "Piña"
$ws.Range("K131").Value = "Caramelo"
$ws.Range("L131").Value = "Primera"
$ws.Range("M131").Value = 60
$ws.Range("N131").Value = 21000
$ws.Range("O131").Value = 22000
$ws.Range("P131").Value = 21500
$ws.Range("Q131").Value = '$/caja 12 unidades'
$ws.Range("R131").Value = "Ecuador"
$ws.Range("S131").Value = 1792
$ws.Range("T131").Value = 12

# New row 132 - Segunda, $/caja 14 unidades
$ws.Range("A132").Value = 7
$ws.Range("B132").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C132").Value = "Ñuble"
$ws.Range("D132").Value = 44461
$ws.Range("E132").Value = 16
$ws.Range("F132").Value = "Fruta"
$ws.Range("G132").Value = 100108
$ws.Range("H132").Value = "Tropicales y subtropicales"
$ws.Range("I132").Value = 100108005
$ws.Range("J132").Value = "Piña"
$ws.Range("K132").Value = "Caramelo"
$ws.Range("L132").Value = "Segunda"
$ws.Range("M132").Value = 60
$ws.Range("N132").Value = 21000
$ws.Range("O132").Value = 22000
$ws.Range("P132").Value = 21500
$ws.Range("Q132").Value = '$/caja 14 unidades'
$ws.Range("R132").Value = "Ecuador"
$ws.Range("S132").Value = 1536
$ws.Range("T132").Value = 14
